# Generate Report for Handoff
#
# Updates the localization-status report after a new handoff-xliff
# generation run:
#   - "Overview" sheet: refresh the "Latest HO Xliff Generate Date" (col G)
#     for the files that just had handoff files generated.
#   - "zh-cn" / "de-de" sheets: refresh their own "Latest Handoff Datetime"
#     (col H) for the same rows, and set "Priority" (col E) to "ht" for
#     those rows (the handoff type used for this run).

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 14)

# --- Overview sheet: bump the handoff-generation timestamp -----------------
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-16 14:21:38"
}

# --- zh-cn sheet: new handoff datetime + priority ---------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("H$r").Value = "2016-08-16 14:21:31"
    $wsZhCn.Range("E$r").Value = "ht"
}

# --- de-de sheet: new handoff datetime + priority ---------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("H$r").Value = "2016-08-16 14:21:38"
    $wsDeDe.Range("E$r").Value = "ht"
}
